# CSci130_Project_gradingsheet.xlsx - "DB Nonsense Taken Care Of"
#
# The grading sheet already has an "x" marker (shared string) in column D
# for every row that was reviewed, except rows 21-23 which were missed.
# Fill those in to match the rest of the sheet, then leave the view/
# selection where the author left it (scrolled back up near row 10,
# with D10 as the active cell) instead of down at D34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSci130_ProjectFall2018")

$ws.Range("D21").Value = "x"
$ws.Range("D22").Value = "x"
$ws.Range("D23").Value = "x"

# Restore the view to scrolled-up-near-the-top with D10 selected (was
# scrolled down to A16/D34 before).
$ws.Range("D10").Select() | Out-Null
